$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.41%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.12%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.158"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07603"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.85%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.326"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.23%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.664"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.41%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9267"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.92%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.13%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1198"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.31%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.36%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09015"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.73%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04143"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.85%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1054"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.26%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001288"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.48%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005805"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.29%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.46%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3355"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.59%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.580"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.58%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1351"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.38%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2806"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.01%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04036"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.71%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001272"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.72%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004065"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.47%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001271"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.15%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02419"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.92%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05169"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.17%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007726"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1300"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.04%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007610"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "16.01%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003303"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "72.52%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008578"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.87%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.43%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006590"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.99%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2686"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-38.82%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004204"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "2.54%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.14%"
